$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns for every coin row with the
# latest scrape. Every Price cell in this sheet is stored as TEXT (prices are
# dotted/grouped strings like "64.722.57" that are not valid numbers, and even
# the plain-decimal ones such as "1.00" or "0.490" must keep their significant
# trailing zeros), so each Price cell is forced to the "@" text format before the
# write and then restored to the default "Normal" style so no stray number format
# is left on the cell - this avoids Excel silently re-interpreting e.g. "1.00" as
# the number 1.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.722.57'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.116.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -7.87%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.06%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -3.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.116.41'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.122'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.07%  '
$ws.Range('E11').Value = '  -5.77%  '
$ws.Range('E12').Value = '  -5.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.664.61'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.83%  '
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.47'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -8.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.640.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000160'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.123.16'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -8.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.66'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.63'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '353.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.16'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.66'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.490'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.267.66'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000113'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -8.30%  '
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('E32').Value = '  -4.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.65'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.20'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.27%  '
$ws.Range('E35').Value = '  -6.20%  '
$ws.Range('E36').Value = '  -5.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.02'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.02%  '
$ws.Range('E38').Value = '  -6.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.823'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.01'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.77%  '
$ws.Range('E41').Value = '  -1.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.624.72'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.73%  '
$ws.Range('E43').Value = '  -3.08%  '
$ws.Range('E44').Value = '  -7.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.25'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('E47').Value = '  -3.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.63'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '316.55'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.67%  '
$ws.Range('E50').Value = '  -4.98%  '
